# Update "想去人数" (column F) counts on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 1578
    3  = 52
    4  = 1034
    5  = 32
    7  = 2707
    9  = 1737
    11 = 77
    12 = 589
    13 = 29
    14 = 17
    15 = 113
    17 = 83
    18 = 16
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
